# Applies the cryptos price/volume refresh described in the commit diff.
# Only the cells that actually changed are touched; D-column numeric-looking
# text values are forced to stay text (matching the original inlineStr type)
# by temporarily applying a text number format and then restoring the style.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$changes = @(
    @{ Row = 2; Col = 'D'; Value = '43.201.62' }
    @{ Row = 2; Col = 'E'; Value = '  -4.99%  ' }
    @{ Row = 3; Col = 'D'; Value = '2.233.22' }
    @{ Row = 3; Col = 'E'; Value = '  -5.73%  ' }
    @{ Row = 4; Col = 'E'; Value = '  +0.19%  ' }
    @{ Row = 5; Col = 'D'; Value = '319.31' }
    @{ Row = 5; Col = 'E'; Value = '  +2.63%  ' }
    @{ Row = 6; Col = 'D'; Value = '99.65' }
    @{ Row = 6; Col = 'E'; Value = '  -9.84%  ' }
    @{ Row = 7; Col = 'E'; Value = '  -7.85%  ' }
    @{ Row = 8; Col = 'E'; Value = '  +0.07%  ' }
    @{ Row = 9; Col = 'D'; Value = '0.567' }
    @{ Row = 10; Col = 'D'; Value = '36.69' }
    @{ Row = 10; Col = 'E'; Value = '  -11.40%  ' }
    @{ Row = 11; Col = 'E'; Value = '  -2.19%  ' }
    @{ Row = 12; Col = 'D'; Value = '0.0822' }
    @{ Row = 12; Col = 'E'; Value = '  -10.75%  ' }
    @{ Row = 13; Col = 'D'; Value = '7.72' }
    @{ Row = 13; Col = 'E'; Value = '  -9.35%  ' }
    @{ Row = 14; Col = 'E'; Value = '  -3.54%  ' }
    @{ Row = 15; Col = 'B'; Value = 'Polygon' }
    @{ Row = 15; Col = 'C'; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' }
    @{ Row = 15; Col = 'D'; Value = '0.867' }
    @{ Row = 15; Col = 'E'; Value = '  -12.25%  ' }
    @{ Row = 16; Col = 'B'; Value = 'WrappedliquidstakedEther2.0' }
    @{ Row = 16; Col = 'C'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth' }
    @{ Row = 16; Col = 'D'; Value = '2.573.06' }
    @{ Row = 16; Col = 'E'; Value = '  -5.60%  ' }
    @{ Row = 17; Col = 'D'; Value = '14.13' }
    @{ Row = 17; Col = 'E'; Value = '  -8.43%  ' }
    @{ Row = 18; Col = 'D'; Value = '2.220.19' }
    @{ Row = 18; Col = 'E'; Value = '  -6.15%  ' }
    @{ Row = 19; Col = 'D'; Value = '43.038.04' }
    @{ Row = 19; Col = 'E'; Value = '  -5.19%  ' }
    @{ Row = 20; Col = 'D'; Value = '14.55' }
    @{ Row = 20; Col = 'E'; Value = '  +4.49%  ' }
    @{ Row = 21; Col = 'D'; Value = '0.0₃0970' }
    @{ Row = 21; Col = 'E'; Value = '  -9.30%  ' }
    @{ Row = 22; Col = 'D'; Value = '6.51' }
    @{ Row = 22; Col = 'E'; Value = '  -11.29%  ' }
    @{ Row = 23; Col = 'D'; Value = '65.43' }
    @{ Row = 23; Col = 'E'; Value = '  -10.99%  ' }
    @{ Row = 24; Col = 'D'; Value = '3.18' }
    @{ Row = 24; Col = 'E'; Value = '  -8.30%  ' }
    @{ Row = 25; Col = 'D'; Value = '236.10' }
    @{ Row = 25; Col = 'E'; Value = '  -9.10%  ' }
    @{ Row = 26; Col = 'D'; Value = '2.16' }
    @{ Row = 26; Col = 'E'; Value = '  -6.70%  ' }
    @{ Row = 27; Col = 'D'; Value = '1.00' }
    @{ Row = 27; Col = 'E'; Value = '  +0.24%  ' }
    @{ Row = 28; Col = 'D'; Value = '10.21' }
    @{ Row = 28; Col = 'E'; Value = '  -8.62%  ' }
    @{ Row = 29; Col = 'E'; Value = '  -7.13%  ' }
    @{ Row = 30; Col = 'E'; Value = '  -13.67%  ' }
    @{ Row = 31; Col = 'D'; Value = '0.0888' }
    @{ Row = 31; Col = 'E'; Value = '  -8.48%  ' }
    @{ Row = 32; Col = 'D'; Value = '20.58' }
    @{ Row = 32; Col = 'E'; Value = '  -8.31%  ' }
    @{ Row = 33; Col = 'D'; Value = '157.78' }
    @{ Row = 33; Col = 'E'; Value = '  -7.52%  ' }
    @{ Row = 34; Col = 'D'; Value = '33.92' }
    @{ Row = 34; Col = 'E'; Value = '  -11.32%  ' }
    @{ Row = 35; Col = 'E'; Value = '  -5.25%  ' }
    @{ Row = 36; Col = 'D'; Value = '3.34' }
    @{ Row = 36; Col = 'E'; Value = '  +12.34%  ' }
    @{ Row = 37; Col = 'D'; Value = '2.00' }
    @{ Row = 37; Col = 'E'; Value = '  +14.63%  ' }
    @{ Row = 38; Col = 'E'; Value = '  -6.59%  ' }
    @{ Row = 39; Col = 'D'; Value = '4.50' }
    @{ Row = 39; Col = 'E'; Value = '  -7.54%  ' }
    @{ Row = 40; Col = 'E'; Value = '  -8.62%  ' }
    @{ Row = 41; Col = 'D'; Value = '3.60' }
    @{ Row = 41; Col = 'E'; Value = '  -9.28%  ' }
    @{ Row = 42; Col = 'D'; Value = '0.0323' }
    @{ Row = 42; Col = 'E'; Value = '  -9.91%  ' }
    @{ Row = 43; Col = 'D'; Value = '1.00' }
    @{ Row = 43; Col = 'E'; Value = '  +0.21%  ' }
    @{ Row = 44; Col = 'D'; Value = '1.825.12' }
    @{ Row = 44; Col = 'E'; Value = '  +9.51%  ' }
    @{ Row = 45; Col = 'D'; Value = '12.12' }
    @{ Row = 45; Col = 'E'; Value = '  -6.02%  ' }
    @{ Row = 46; Col = 'D'; Value = '88.21' }
    @{ Row = 46; Col = 'E'; Value = '  -11.40%  ' }
    @{ Row = 47; Col = 'B'; Value = 'Algorand' }
    @{ Row = 47; Col = 'C'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' }
    @{ Row = 47; Col = 'D'; Value = '0.208' }
    @{ Row = 47; Col = 'E'; Value = '  -11.18%  ' }
    @{ Row = 48; Col = 'D'; Value = '5.50' }
    @{ Row = 48; Col = 'E'; Value = '  -0.09%  ' }
    @{ Row = 49; Col = 'B'; Value = 'ordi' }
    @{ Row = 49; Col = 'C'; Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi' }
    @{ Row = 49; Col = 'D'; Value = '77.92' }
    @{ Row = 49; Col = 'E'; Value = '  -6.28%  ' }
    @{ Row = 50; Col = 'B'; Value = 'MultiversX' }
    @{ Row = 50; Col = 'C'; Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld' }
    @{ Row = 50; Col = 'D'; Value = '60.73' }
    @{ Row = 50; Col = 'E'; Value = '  -13.43%  ' }
    @{ Row = 51; Col = 'B'; Value = 'FraxShare' }
    @{ Row = 51; Col = 'C'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' }
    @{ Row = 51; Col = 'D'; Value = '8.57' }
    @{ Row = 51; Col = 'E'; Value = '  -7.39%  ' }
)

foreach ($change in $changes) {
    $colIndex = switch ($change.Col) {
        'A' { 1 }
        'B' { 2 }
        'C' { 3 }
        'D' { 4 }
        'E' { 5 }
    }
    $cell = $ws.Cells.Item($change.Row, $colIndex)
    if ($change.Col -eq 'D') {
        # Force text type so numeric-looking prices (e.g. '1.00') are not
        # reinterpreted as numbers, then drop back to the default cell style
        # so no stray number-format style gets attached to the cell.
        $cell.NumberFormat = '@'
        $cell.Value = $change.Value
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $change.Value
    }
}
